$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 52, shifting existing rows 52-102 down to 53-103
$ws.Rows.Item(52).Insert()

# Populate the new row 52 with the new record
$ws.Cells.Item(52, 1).Value = 1
$ws.Cells.Item(52, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(52, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(52, 4).Value = 44827
$ws.Cells.Item(52, 5).Value = 15
$ws.Cells.Item(52, 6).Value = "Fruta"
$ws.Cells.Item(52, 7).Value = 100102
$ws.Cells.Item(52, 8).Value = "Cítricos"
$ws.Cells.Item(52, 9).Value = 100102005
$ws.Cells.Item(52, 10).Value = "Naranja"
$ws.Cells.Item(52, 11).Value = "Lane Late"
$ws.Cells.Item(52, 12).Value = "Tercera"
$ws.Cells.Item(52, 13).Value = 300
$ws.Cells.Item(52, 14).Value = 500
$ws.Cells.Item(52, 15).Value = 600
$ws.Cells.Item(52, 16).Value = 550
$ws.Cells.Item(52, 17).Value = "`$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(52, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(52, 19).Value = 550
$ws.Cells.Item(52, 20).Value = 1
